$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/11_tokiko"
$ws.Range("B2").Value = "pngimages/11_compass.png"
$ws.Range("C2").Value = "trainingimages/27_pakapa"
$ws.Range("D2").Value = "pngimages/27_kiwi.png"
$ws.Range("E2").Value = -0.5
$ws.Range("F2").Value = 0.5

# Row 3
$ws.Range("A3").Value = "trainingimages/17_kotako"
$ws.Range("B3").Value = "pngimages/17_cracker.png"
$ws.Range("C3").Value = "trainingimages/25_tapapi"
$ws.Range("D3").Value = "pngimages/25_apple.png"
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = -0.5

# Row 4
$ws.Range("A4").Value = "trainingimages/12_pokika"
$ws.Range("B4").Value = "pngimages/12_pie.png"
$ws.Range("C4").Value = "trainingimages/07_pitapi"
$ws.Range("D4").Value = "pngimages/07_suitcase.png"
